$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "48÷2=24, 0"
$t.Cell(1,2).Range.Text = "23÷8=2, 7"
$t.Cell(1,3).Range.Text = "37÷9=4, 1"
$t.Cell(1,4).Range.Text = "18÷4=4, 2"
$t.Cell(1,5).Range.Text = "43÷6=7, 1"
$t.Cell(5,1).Range.Text = "36÷7=5, 1"
$t.Cell(5,2).Range.Text = "37÷3=12, 1"
$t.Cell(5,3).Range.Text = "35÷3=11, 2"
$t.Cell(5,4).Range.Text = "40÷6=6, 4"
$t.Cell(5,5).Range.Text = "79÷4=19, 3"
$t.Cell(9,1).Range.Text = "69÷8=8, 5"
$t.Cell(9,2).Range.Text = "79÷6=13, 1"
$t.Cell(9,3).Range.Text = "59÷5=11, 4"
$t.Cell(9,4).Range.Text = "99÷2=49, 1"
$t.Cell(9,5).Range.Text = "82÷6=13, 4"
$t.Cell(13,1).Range.Text = "60÷3=20, 0"
$t.Cell(13,2).Range.Text = "81÷6=13, 3"
$t.Cell(13,3).Range.Text = "42÷6=7, 0"
$t.Cell(13,4).Range.Text = "62÷6=10, 2"
$t.Cell(13,5).Range.Text = "17÷3=5, 2"
$t.Cell(17,1).Range.Text = "75÷6=12, 3"
$t.Cell(17,2).Range.Text = "84÷2=42, 0"
$t.Cell(17,3).Range.Text = "89÷4=22, 1"
$t.Cell(17,4).Range.Text = "20÷4=5, 0"
$t.Cell(17,5).Range.Text = "11÷5=2, 1"
